$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.347
$ws.Range("C4").Value = -11.699
$ws.Range("B6").Value = 5.848999999999999
$ws.Range("B7").Value = 5.733
$ws.Range("D7").Value = -7.536
$ws.Range("D8").Value = -8.122000000000002
$ws.Range("C9").Value = -10.827
$ws.Range("D10").Value = -7.593000000000001
$ws.Range("C12").Value = -10.956
$ws.Range("D13").Value = -7.772
$ws.Range("B16").Value = 5.324
$ws.Range("D16").Value = -8.052
$ws.Range("C17").Value = -13.474
$ws.Range("C18").Value = -12.588
$ws.Range("C19").Value = -11.949
$ws.Range("B20").Value = 8.914999999999999
$ws.Range("C20").Value = -12.243
$ws.Range("C26").Value = -12.635
$ws.Range("B28").Value = 5.196000000000001
$ws.Range("B29").Value = 5.213
$ws.Range("D30").Value = -7.289999999999999
$ws.Range("C31").Value = -13.298
$ws.Range("B32").Value = 7.084999999999999
$ws.Range("C39").Value = -12.278
$ws.Range("B40").Value = 9.236000000000001
$ws.Range("C40").Value = -12.09
$ws.Range("D40").Value = -8.43
$ws.Range("C41").Value = -12.02
$ws.Range("C42").Value = -12.364
$ws.Range("C43").Value = -12.182
$ws.Range("D44").Value = -7.513000000000001
$ws.Range("B46").Value = 5.447000000000001
$ws.Range("C47").Value = -13.123
$ws.Range("C48").Value = -11.844
$ws.Range("B51").Value = 5.399
$ws.Range("B52").Value = 5.411
$ws.Range("B57").Value = 5.263
$ws.Range("B59").Value = 4.878000000000001
$ws.Range("B62").Value = 5.494
$ws.Range("C63").Value = -10.956
$ws.Range("C64").Value = -11.151
$ws.Range("B66").Value = 5.511
$ws.Range("B73").Value = 5.884
$ws.Range("B74").Value = 9.186999999999999
$ws.Range("C76").Value = -12.288
$ws.Range("C81").Value = -13.148
$ws.Range("C89").Value = -13.278
$ws.Range("D89").Value = -8.32
$ws.Range("D91").Value = -7.486
$ws.Range("B92").Value = 4.891
$ws.Range("C94").Value = -11.828
$ws.Range("B100").Value = 6.026999999999999
